# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) cells for the existing rows, and also
# fixes row order for ARBITRUM/Aave (rows 50-51 swap which coin occupies which
# row, along with their refreshed price/volume).
#
# Several new Price values happen to parse as plain numbers (e.g. "251.59"),
# but the source column stores prices as text (note other rows like
# "43.118.09" with two dots, which can't be numbers). Assigning such a
# string via .Value would make Excel auto-convert the cell to a number. To
# keep those cells as text - matching the original file - we prefix the
# value with a leading apostrophe (the classic "force text" trick) and then
# explicitly restore the cell's style to "Normal" afterwards so no stray
# quote-prefix/number-format styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.924.84'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.285.98'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''251.59'
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").Value = '''0.640'
$ws.Range("E6").Value = '  +2.36%  '
$ws.Range("D7").Value = '''73.39'
$ws.Range("E7").Value = '  +2.10%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '''0.634'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").Value = '''38.96'
$ws.Range("E10").Value = '  -5.48%  '
$ws.Range("D11").Value = '''0.0976'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '''59.09'
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").Value = '''7.42'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").Value = '2.631.39'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("D16").Value = '''15.20'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").Value = '''0.871'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '2.280.85'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = '42.856.46'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").Value = '0.0₃0999'
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D21").Value = '''6.28'
$ws.Range("E21").Value = '  +0.35%  '
$ws.Range("D22").Value = '''72.51'
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").Value = '''237.22'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").Value = '''2.21'
$ws.Range("E24").Value = '  +3.96%  '
$ws.Range("E25").Value = '  -2.56%  '
$ws.Range("D26").Value = '''11.56'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '''2.41'
$ws.Range("E28").Value = '  -2.06%  '
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").Value = '''2.19'
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("D31").Value = '''166.46'
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").Value = '''21.01'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '''6.45'
$ws.Range("E33").Value = '  +3.98%  '
$ws.Range("D34").Value = '''0.126'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").Value = '''0.0824'
$ws.Range("E35").Value = '  +4.97%  '
$ws.Range("D36").Value = '''30.87'
$ws.Range("E36").Value = '  +9.89%  '
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("D38").Value = '''4.58'
$ws.Range("E38").Value = '  +10.38%  '
$ws.Range("D39").Value = '''4.75'
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("D40").Value = '''0.0308'
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("D41").Value = '''14.14'
$ws.Range("E41").Value = '  +13.79%  '
$ws.Range("D42").Value = '''2.33'
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("D43").Value = '''5.92'
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  +6.66%  '
$ws.Range("D45").Value = '''9.15'
$ws.Range("E45").Value = '  +2.65%  '
$ws.Range("D46").Value = '''61.59'
$ws.Range("E46").Value = '  -5.13%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '''100.84'
$ws.Range("E50").Value = '  +6.66%  '
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").Value = '''1.17'
$ws.Range("E51").Value = '  -2.74%  '

# Reset style on cells that required the text-forcing apostrophe prefix,
# so no stray quotePrefix / number-format style survives on the cell.
foreach ($addr in @("D5","D6","D7","D9","D10","D11","D12","D13","D16","D17","D21","D22","D23","D24","D26","D28","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D45","D46","D50","D51")) {
    $ws.Range($addr).Style = "Normal"
}
